$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in win/loss/tie record for every data row (2 through 52)
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 94   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 68   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
